$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RYI")
$ws.Columns("R").Insert()
